# Power_BusInfo.xlsx: bump template version and rename a few field-key
# labels on the "scenarioA" sheet.
#   C2  "v0.0.3r"   -> "v0.0.4r"
#   A4  "Excl."     -> "excl"
#   K4  "comYear"   -> "YearCom"
#   L4  "decomYear" -> "YearDecom"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = "v0.0.4r"
$ws.Range("A4").Value = "excl"
$ws.Range("K4").Value = "YearCom"
$ws.Range("L4").Value = "YearDecom"
